$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.991.74'
$ws.Range("E2").Value = '  +1.95%  '
$ws.Range("D3").Value = '1.908.83'
$ws.Range("E3").Value = '  +2.26%  '
$ws.Range("E4").Value = '  -0.79%  '
$ws.Range("D5").Value = "'315.27"
$ws.Range("E5").Value = '  +1.09%  '
$ws.Range("E6").Value = '  -0.84%  '
$ws.Range("D7").Value = "'0.4802"
$ws.Range("E7").Value = '  +0.38%  '
$ws.Range("D8").Value = "'0.3801"
$ws.Range("E8").Value = '  +0.93%  '
$ws.Range("D9").Value = "'0.07361"
$ws.Range("E9").Value = '  +0.33%  '
$ws.Range("D10").Value = "'0.9328"
$ws.Range("E10").Value = '  -0.32%  '
$ws.Range("D11").Value = "'20.80"
$ws.Range("E11").Value = '  +0.47%  '
$ws.Range("D12").Value = "'0.07765"
$ws.Range("E12").Value = '  -0.86%  '
$ws.Range("D13").Value = '1.853.71'
$ws.Range("E13").Value = '  -1.12%  '
$ws.Range("E14").Value = '  +1.04%  '
$ws.Range("E15").Value = '  +1.25%  '
$ws.Range("D16").Value = "'91.76"
$ws.Range("E16").Value = '  +1.43%  '
$ws.Range("D17").Value = "'1.005"
$ws.Range("E17").Value = '  -0.79%  '
$ws.Range("D18").Value = "'0.000008834"
$ws.Range("E18").Value = '  -0.74%  '
$ws.Range("D19").Value = "'1.003"
$ws.Range("D20").Value = '28.022.28'
$ws.Range("E20").Value = '  +1.81%  '
$ws.Range("E21").Value = '  +0.31%  '
$ws.Range("D22").Value = "'5.166"
$ws.Range("E22").Value = '  +0.95%  '
$ws.Range("D23").Value = '2.162.14'
$ws.Range("E23").Value = '  +1.49%  '
$ws.Range("E24").Value = '  +1.65%  '
$ws.Range("D25").Value = "'155.61"
$ws.Range("E25").Value = '  +0.41%  '
$ws.Range("D26").Value = "'1.914"
$ws.Range("E26").Value = '  -1.20%  '
$ws.Range("E27").Value = '  -0.04%  '
$ws.Range("D28").Value = "'2.128"
$ws.Range("E28").Value = '  +5.22%  '
$ws.Range("D29").Value = "'116.86"
$ws.Range("D30").Value = "'4.954"
$ws.Range("E30").Value = '  -0.43%  '
$ws.Range("D31").Value = "'0.08934"
$ws.Range("E31").Value = '  +0.40%  '
$ws.Range("D32").Value = "'3.295"
$ws.Range("E32").Value = '  -1.00%  '
$ws.Range("E33").Value = '  +3.52%  '
$ws.Range("D34").Value = "'0.7737"
$ws.Range("E34").Value = '  +2.28%  '
$ws.Range("D35").Value = "'4.678"
$ws.Range("E35").Value = '  +1.63%  '
$ws.Range("D36").Value = "'2.617"
$ws.Range("D37").Value = "'0.02053"
$ws.Range("E37").Value = '  +0.73%  '
$ws.Range("D38").Value = "'1.112"
$ws.Range("E38").Value = '  -0.68%  '
$ws.Range("D39").Value = "'0.05305"
$ws.Range("E39").Value = '  +0.83%  '
$ws.Range("D40").Value = "'3.004"
$ws.Range("E40").Value = '  +0.52%  '
$ws.Range("D41").Value = "'0.5489"
$ws.Range("E41").Value = '  +3.14%  '
$ws.Range("D42").Value = "'7.020"
$ws.Range("E42").Value = '  -0.93%  '
$ws.Range("D43").Value = "'0.1525"
$ws.Range("E43").Value = '  -0.01%  '
$ws.Range("D44").Value = "'8.464"
$ws.Range("E44").Value = '  -1.29%  '
$ws.Range("D45").Value = "'10.64"
$ws.Range("E45").Value = '  -0.08%  '
$ws.Range("E46").Value = '  +0.43%  '
$ws.Range("D47").Value = "'108.17"
$ws.Range("E47").Value = '  +5.18%  '
$ws.Range("D48").Value = "'1.004"
$ws.Range("E48").Value = '  -0.89%  '
$ws.Range("E49").Value = '  -0.44%  '
$ws.Range("D50").Value = "'67.88"
$ws.Range("E50").Value = '  +0.80%  '
$ws.Range("D51").Value = "'0.06075"
$ws.Range("E51").Value = '  -0.05%  '
